# Weekly fruit/vegetable update: insert two new rows of data (the new
# week's prices) above the existing block starting at row 278. This
# pushes the former rows 278-283 down to rows 280-285 and adds two new
# rows (278-279) with the new observations dated 2021-09-09 (serial 44448).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 278 - this shifts the old 278-283 block
# down to 280-285 and (per Excel's normal insert behavior) copies the
# number format (date style) of column D down from the row below.
$ws.Rows.Item(278).Insert()
$ws.Rows.Item(278).Insert()

# --- New row 278 ---
$ws.Cells.Item(278, 1).Value = 4
$ws.Cells.Item(278, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(278, 3).Value = "Los Lagos"
$ws.Cells.Item(278, 4).Value2 = 44448
$ws.Cells.Item(278, 5).Value = 10
$ws.Cells.Item(278, 6).Value = 100112004
$ws.Cells.Item(278, 7).Value = "Cebolla"
$ws.Cells.Item(278, 8).Value = "Morada(o)"
$ws.Cells.Item(278, 9).Value = "1a (guarda)"
$ws.Cells.Item(278, 10).Value = 140
$ws.Cells.Item(278, 11).Value = 13500
$ws.Cells.Item(278, 12).Value = 14000
$ws.Cells.Item(278, 13).Value = 13750
$ws.Cells.Item(278, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item(278, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(278, 16).Value = 764
$ws.Cells.Item(278, 17).Value = 18
$ws.Cells.Item(278, 18).Value = "Hortaliza"

# --- New row 279 ---
$ws.Cells.Item(279, 1).Value = 4
$ws.Cells.Item(279, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(279, 3).Value = "Los Lagos"
$ws.Cells.Item(279, 4).Value2 = 44448
$ws.Cells.Item(279, 5).Value = 10
$ws.Cells.Item(279, 6).Value = 100112004
$ws.Cells.Item(279, 7).Value = "Cebolla"
$ws.Cells.Item(279, 8).Value = "Sin especificar"
$ws.Cells.Item(279, 9).Value = "1a (guarda)"
$ws.Cells.Item(279, 10).Value = 450
$ws.Cells.Item(279, 11).Value = 7000
$ws.Cells.Item(279, 12).Value = 7000
$ws.Cells.Item(279, 13).Value = 7000
$ws.Cells.Item(279, 14).Value = '$/malla 16 kilos'
$ws.Cells.Item(279, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(279, 16).Value = 438
$ws.Cells.Item(279, 17).Value = 16
$ws.Cells.Item(279, 18).Value = "Hortaliza"
